# "La matriz debe ser potencia de 2"
# Add a new sheet "Hoja3" (placed after "Hoja2") containing a fresh copy of
# the original 8x8 adjacency matrix (Hoja1!A1:H8) so it can be squared again
# as a power-of-two sized matrix. Also update the previously-active sheet
# (Hoja2)'s saved selection/scroll state, since it is no longer the active
# tab once Hoja3 becomes selected.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# Insert the new worksheet right after "Hoja2" and name it "Hoja3".
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Hoja3"

# Copy the 8x8 matrix values from Hoja1 (A1:H8) into the new sheet.
$ws1.Range("A1:H8").Copy()
$ws3.Range("A1").PasteSpecial()
$ws3.Range("A1:H8").NumberFormat = "0"

# Restore Hoja2's selection to its last-used cell (no longer the active tab).
$ws2.Range("G12").Select()

# Hoja3 becomes the active sheet/tab with its own selection state.
$ws3.Range("I1:T30").Select()
